$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E17").Value = 20.9
$ws.Range("E18").Value = 14.4
$ws.Range("E19").Value = 12.6
$ws.Range("E20").Value = 9

$ws.Range("R12").Select()
